$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 204 (Especial/Primera/Segunda @ O'Higgins,
# dated 44399) so the old rows 204-209 shift down to become rows 206-211, matching the
# diff's unchanged tail content.
$ws.Range("A204:T205").Insert()

# New row 204: Kiwi Hayward Primera, $/bandeja 10 kilos, Region Metropolitana, date 44448
$ws.Range("A204").Value = 8
$ws.Range("B204").Value = "Terminal La Palmera de La Serena"
$ws.Range("C204").Value = "Coquimbo"
$ws.Range("D204").Value = 44448
$ws.Range("E204").Value = 4
$ws.Range("F204").Value = "Fruta"
$ws.Range("G204").Value = 100101
$ws.Range("H204").Value = "Berries"
$ws.Range("I204").Value = 100101007
$ws.Range("J204").Value = "Kiwi"
$ws.Range("K204").Value = "Hayward"
$ws.Range("L204").Value = "Primera"
$ws.Range("M204").Value = 300
$ws.Range("N204").Value = 9000
$ws.Range("O204").Value = 10000
$ws.Range("P204").Value = 9500
$ws.Range("Q204").Value = "$/bandeja 10 kilos"
$ws.Range("R204").Value = "Región Metropolitana"
$ws.Range("S204").Value = 950
$ws.Range("T204").Value = 10

# New row 205: Kiwi Hayward Segunda, $/bandeja 10 kilos, Region Metropolitana, date 44448
$ws.Range("A205").Value = 8
$ws.Range("B205").Value = "Terminal La Palmera de La Serena"
$ws.Range("C205").Value = "Coquimbo"
$ws.Range("D205").Value = 44448
$ws.Range("E205").Value = 4
$ws.Range("F205").Value = "Fruta"
$ws.Range("G205").Value = 100101
$ws.Range("H205").Value = "Berries"
$ws.Range("I205").Value = 100101007
$ws.Range("J205").Value = "Kiwi"
$ws.Range("K205").Value = "Hayward"
$ws.Range("L205").Value = "Segunda"
$ws.Range("M205").Value = 300
$ws.Range("N205").Value = 7000
$ws.Range("O205").Value = 8000
$ws.Range("P205").Value = 7500
$ws.Range("Q205").Value = "$/bandeja 10 kilos"
$ws.Range("R205").Value = "Región Metropolitana"
$ws.Range("S205").Value = 750
$ws.Range("T205").Value = 10
